$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row's cells to be entered as literal text (matching the
# existing "inlineStr" cells above them) rather than letting Excel's
# smart-typing turn "2025-01-20" into a date serial or "42.6" into a number.
$ws.Range("A83:B83").NumberFormat = "@"
$ws.Range("A83").Value = "2025-01-20"
$ws.Range("B83").Value = "42.6"

# Reset back to the default (unstyled) cell style so the appended row
# doesn't pick up an explicit style index, matching the other data rows.
$ws.Range("A83:B83").Style = "Normal"
